$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, pushing existing rows 14-43 down to 15-44
$ws.Rows.Item(14).Insert()

# Fill the new row 14 with the new data point (values come from the diff)
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "Vega Monumental Concepción"
$ws.Range("C14").Value = "Bíobío"
$ws.Range("D14").Value = 44883
$ws.Range("E14").Value = 8
$ws.Range("F14").Value = 100114007
$ws.Range("G14").Value = "Jengibre"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 60
$ws.Range("K14").Value = 14000
$ws.Range("L14").Value = 15000
$ws.Range("M14").Value = 14500
$ws.Range("N14").Value = '$/caja 13 kilos'
$ws.Range("O14").Value = "Perú"
$ws.Range("P14").Value = 1115
$ws.Range("Q14").Value = 13
$ws.Range("R14").Value = "Hortaliza"
